#
# Applies the "Welcome" demo-text simplification to the pythoncode sheet:
#   - A2: print("Welcome to the DS-Algo demo");  ->  print("Welcome");
#   - B2: Welcome to the demo                     ->  Welcome
# and leaves the cursor/selection on B8, matching the author's final
# on-screen state when the workbook was saved.
#
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pythoncode")
$ws.Activate()

# Set B2 before A2 so the shared-string table is rebuilt in the same
# order ("Welcome" then the print statement referencing it).
$ws.Range("B2").Value = "Welcome"
$ws.Range("A2").Value = 'print("Welcome");'

$ws.Range("B8").Select()
